$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 189-193: dates "06 08 2020" through "10 08 2020" (column A only)
$ws.Cells.Item(189, 1).Value = "06 08 2020"
$ws.Cells.Item(190, 1).Value = "07 08 2020"
$ws.Cells.Item(191, 1).Value = "08 08 2020"
$ws.Cells.Item(192, 1).Value = "09 08 2020"
$ws.Cells.Item(193, 1).Value = "10 08 2020"

# New column E values for rows 90-94
$ws.Cells.Item(90, 5).Value = 0.11762681762682
$ws.Cells.Item(91, 5).Value = 0.12646103896104
$ws.Cells.Item(92, 5).Value = 0.15536315536316
$ws.Cells.Item(93, 5).Value = 0.19444444444444
$ws.Cells.Item(94, 5).Value = 0.11825396825397

# New data for rows 184-188, columns B:D and F:BE (E stays blank)
# Row 184
$ws.Cells.Item(184, 2).Value = 0.076436333590527
$ws.Cells.Item(184, 3).Value = 0.070493328690589
$ws.Cells.Item(184, 4).Value = 0.062323286146406
$ws.Cells.Item(184, 6).Value = 0.066492995665994
$ws.Cells.Item(184, 7).Value = 0.05482500355099
$ws.Cells.Item(184, 8).Value = 0.066251497359363
$ws.Cells.Item(184, 9).Value = 0.058588400624574
$ws.Cells.Item(184, 10).Value = 0.051838027941349
$ws.Cells.Item(184, 11).Value = 0.056429124609332
$ws.Cells.Item(184, 12).Value = 0.051866411379361
$ws.Cells.Item(184, 13).Value = 0.057522731592719
$ws.Cells.Item(184, 14).Value = 0.065558643370766
$ws.Cells.Item(184, 15).Value = 0.056282133458982
$ws.Cells.Item(184, 16).Value = 0.059644129817532
$ws.Cells.Item(184, 17).Value = 0.067662669923942
$ws.Cells.Item(184, 18).Value = 0.052422683395611
$ws.Cells.Item(184, 19).Value = 0.049648632612831
$ws.Cells.Item(184, 20).Value = 0.053652735861181
$ws.Cells.Item(184, 21).Value = 0.04578304242305
$ws.Cells.Item(184, 22).Value = 0.054816975700815
$ws.Cells.Item(184, 23).Value = 0.056684182278314
$ws.Cells.Item(184, 24).Value = 0.044706949126838
$ws.Cells.Item(184, 25).Value = 0.055555957837512
$ws.Cells.Item(184, 26).Value = 0.053639281068254
$ws.Cells.Item(184, 27).Value = 0.052848493447133
$ws.Cells.Item(184, 28).Value = 0.058096191016408
$ws.Cells.Item(184, 29).Value = 0.077336772276886
$ws.Cells.Item(184, 30).Value = 0.064434466361474
$ws.Cells.Item(184, 31).Value = 0.069779213866266
$ws.Cells.Item(184, 32).Value = 0.056533870850582
$ws.Cells.Item(184, 33).Value = 0.056923978130727
$ws.Cells.Item(184, 34).Value = 0.066353951120719
$ws.Cells.Item(184, 35).Value = 0.061879254192375
$ws.Cells.Item(184, 36).Value = 0.058752794246043
$ws.Cells.Item(184, 37).Value = 0.06653057017163
$ws.Cells.Item(184, 38).Value = 0.060311582778679
$ws.Cells.Item(184, 39).Value = 0.057367102006891
$ws.Cells.Item(184, 40).Value = 0.051755122361175
$ws.Cells.Item(184, 41).Value = 0.058451133695989
$ws.Cells.Item(184, 42).Value = 0.057509943126288
$ws.Cells.Item(184, 43).Value = 0.053658524132969
$ws.Cells.Item(184, 44).Value = 0.062339364043975
$ws.Cells.Item(184, 45).Value = 0.063637440618334
$ws.Cells.Item(184, 46).Value = 0.059764738977201
$ws.Cells.Item(184, 47).Value = 0.065442609227615
$ws.Cells.Item(184, 48).Value = 0.056742558680475
$ws.Cells.Item(184, 49).Value = 0.051330394589238
$ws.Cells.Item(184, 50).Value = 0.06773395237829
$ws.Cells.Item(184, 51).Value = 0.045336236755141
$ws.Cells.Item(184, 52).Value = 0.054998454246246
$ws.Cells.Item(184, 53).Value = 0.051246364060736
$ws.Cells.Item(184, 54).Value = 0.049597124052313
$ws.Cells.Item(184, 55).Value = 0.050935057775126
$ws.Cells.Item(184, 56).Value = 0.044193640163265
$ws.Cells.Item(184, 57).Value = 0.062923088996558
# Row 185
$ws.Cells.Item(185, 2).Value = 0.054124264025565
$ws.Cells.Item(185, 3).Value = 0.055117639020507
$ws.Cells.Item(185, 4).Value = 0.06416650134123
$ws.Cells.Item(185, 6).Value = 0.057644228336931
$ws.Cells.Item(185, 7).Value = 0.065275938155212
$ws.Cells.Item(185, 8).Value = 0.058661642774864
$ws.Cells.Item(185, 9).Value = 0.042972325234784
$ws.Cells.Item(185, 10).Value = 0.040424881716329
$ws.Cells.Item(185, 11).Value = 0.04500729435753
$ws.Cells.Item(185, 12).Value = 0.034979720109028
$ws.Cells.Item(185, 13).Value = 0.041518381883286
$ws.Cells.Item(185, 14).Value = 0.067063724933891
$ws.Cells.Item(185, 15).Value = 0.044252702768139
$ws.Cells.Item(185, 16).Value = 0.05177956408543
$ws.Cells.Item(185, 17).Value = 0.056014335418566
$ws.Cells.Item(185, 18).Value = 0.043553628486967
$ws.Cells.Item(185, 19).Value = 0.045986851465373
$ws.Cells.Item(185, 20).Value = 0.051100065100899
$ws.Cells.Item(185, 21).Value = 0.04984345331077
$ws.Cells.Item(185, 22).Value = 0.058856443844999
$ws.Cells.Item(185, 23).Value = 0.057008981095378
$ws.Cells.Item(185, 24).Value = 0.049084786067812
$ws.Cells.Item(185, 25).Value = 0.056929430724886
$ws.Cells.Item(185, 26).Value = 0.045751329350938
$ws.Cells.Item(185, 27).Value = 0.045655493276351
$ws.Cells.Item(185, 28).Value = 0.049581235581171
$ws.Cells.Item(185, 29).Value = 0.077402129619607
$ws.Cells.Item(185, 30).Value = 0.052065561125204
$ws.Cells.Item(185, 31).Value = 0.058982600077502
$ws.Cells.Item(185, 32).Value = 0.050112075549664
$ws.Cells.Item(185, 33).Value = 0.047722876499328
$ws.Cells.Item(185, 34).Value = 0.051077061641307
$ws.Cells.Item(185, 35).Value = 0.049481433241623
$ws.Cells.Item(185, 36).Value = 0.045414558740715
$ws.Cells.Item(185, 37).Value = 0.060171697246022
$ws.Cells.Item(185, 38).Value = 0.050737271249347
$ws.Cells.Item(185, 39).Value = 0.04969181025952
$ws.Cells.Item(185, 40).Value = 0.044402141840993
$ws.Cells.Item(185, 41).Value = 0.05022527281716
$ws.Cells.Item(185, 42).Value = 0.049090990869217
$ws.Cells.Item(185, 43).Value = 0.0432391283909
$ws.Cells.Item(185, 44).Value = 0.067561273457852
$ws.Cells.Item(185, 45).Value = 0.05034352057283
$ws.Cells.Item(185, 46).Value = 0.048020940399804
$ws.Cells.Item(185, 47).Value = 0.053597411923995
$ws.Cells.Item(185, 48).Value = 0.050071335731822
$ws.Cells.Item(185, 49).Value = 0.045181139811424
$ws.Cells.Item(185, 50).Value = 0.055113897125296
$ws.Cells.Item(185, 51).Value = 0.047192796464783
$ws.Cells.Item(185, 52).Value = 0.07010943429699899
$ws.Cells.Item(185, 53).Value = 0.047172516620228
$ws.Cells.Item(185, 54).Value = 0.050418399330215
$ws.Cells.Item(185, 55).Value = 0.048827890542825
$ws.Cells.Item(185, 56).Value = 0.047466481968798
$ws.Cells.Item(185, 57).Value = 0.062244894556552
# Row 186
$ws.Cells.Item(186, 2).Value = 0.089048400519762
$ws.Cells.Item(186, 3).Value = 0.10023830635236
$ws.Cells.Item(186, 4).Value = 0.10625293179318
$ws.Cells.Item(186, 6).Value = 0.08029015016720301
$ws.Cells.Item(186, 7).Value = 0.081057037873541
$ws.Cells.Item(186, 8).Value = 0.084315299337547
$ws.Cells.Item(186, 9).Value = 0.070642558949588
$ws.Cells.Item(186, 10).Value = 0.057591003008516
$ws.Cells.Item(186, 11).Value = 0.067459719126278
$ws.Cells.Item(186, 12).Value = 0.07124749978988699
$ws.Cells.Item(186, 13).Value = 0.088286416294619
$ws.Cells.Item(186, 14).Value = 0.12414696287448
$ws.Cells.Item(186, 15).Value = 0.08803827492068
$ws.Cells.Item(186, 16).Value = 0.082585418403836
$ws.Cells.Item(186, 17).Value = 0.088748150331162
$ws.Cells.Item(186, 18).Value = 0.06391911438995
$ws.Cells.Item(186, 19).Value = 0.073620848194757
$ws.Cells.Item(186, 20).Value = 0.08286683210821399
$ws.Cells.Item(186, 21).Value = 0.073481242882493
$ws.Cells.Item(186, 22).Value = 0.08609908282230801
$ws.Cells.Item(186, 23).Value = 0.062798382855239
$ws.Cells.Item(186, 24).Value = 0.05648004557771
$ws.Cells.Item(186, 25).Value = 0.066361203619201
$ws.Cells.Item(186, 26).Value = 0.06784263115787401
$ws.Cells.Item(186, 27).Value = 0.07110415132304999
$ws.Cells.Item(186, 28).Value = 0.08647154737886099
$ws.Cells.Item(186, 29).Value = 0.1327346700609
$ws.Cells.Item(186, 30).Value = 0.10896298545042
$ws.Cells.Item(186, 31).Value = 0.096842600897769
$ws.Cells.Item(186, 32).Value = 0.07993994794760501
$ws.Cells.Item(186, 33).Value = 0.092296975196289
$ws.Cells.Item(186, 34).Value = 0.09832289672833899
$ws.Cells.Item(186, 35).Value = 0.07546998716511
$ws.Cells.Item(186, 36).Value = 0.069909881292365
$ws.Cells.Item(186, 37).Value = 0.082158723268485
$ws.Cells.Item(186, 38).Value = 0.074761910047597
$ws.Cells.Item(186, 39).Value = 0.06564103653282399
$ws.Cells.Item(186, 40).Value = 0.074868108311404
$ws.Cells.Item(186, 41).Value = 0.097835152275963
$ws.Cells.Item(186, 42).Value = 0.078970738342328
$ws.Cells.Item(186, 43).Value = 0.067591948040383
$ws.Cells.Item(186, 44).Value = 0.079277206460109
$ws.Cells.Item(186, 45).Value = 0.070262393112221
$ws.Cells.Item(186, 46).Value = 0.084838998871181
$ws.Cells.Item(186, 47).Value = 0.09603257035762899
$ws.Cells.Item(186, 48).Value = 0.088476768395417
$ws.Cells.Item(186, 49).Value = 0.079920414610419
$ws.Cells.Item(186, 50).Value = 0.089209316381511
$ws.Cells.Item(186, 51).Value = 0.067299185717288
$ws.Cells.Item(186, 52).Value = 0.073844559627635
$ws.Cells.Item(186, 53).Value = 0.064237170558158
$ws.Cells.Item(186, 54).Value = 0.068389337678137
$ws.Cells.Item(186, 55).Value = 0.07055682613604999
$ws.Cells.Item(186, 56).Value = 0.06868528957787499
$ws.Cells.Item(186, 57).Value = 0.097797369138602
# Row 187
$ws.Cells.Item(187, 2).Value = 0.08755888006453
$ws.Cells.Item(187, 3).Value = 0.1055595803928
$ws.Cells.Item(187, 4).Value = 0.10439094044686
$ws.Cells.Item(187, 6).Value = 0.07621385458515
$ws.Cells.Item(187, 7).Value = 0.07748634443046799
$ws.Cells.Item(187, 8).Value = 0.086305759110459
$ws.Cells.Item(187, 9).Value = 0.069479822704674
$ws.Cells.Item(187, 10).Value = 0.049780718305272
$ws.Cells.Item(187, 11).Value = 0.057775368483819
$ws.Cells.Item(187, 12).Value = 0.078886864504447
$ws.Cells.Item(187, 13).Value = 0.09508879712222
$ws.Cells.Item(187, 14).Value = 0.079508872651957
$ws.Cells.Item(187, 15).Value = 0.09238993318100799
$ws.Cells.Item(187, 16).Value = 0.10017278469385
$ws.Cells.Item(187, 17).Value = 0.098554146464944
$ws.Cells.Item(187, 18).Value = 0.077716570039999
$ws.Cells.Item(187, 19).Value = 0.090767924307337
$ws.Cells.Item(187, 20).Value = 0.09732876637231801
$ws.Cells.Item(187, 21).Value = 0.08643377586497999
$ws.Cells.Item(187, 22).Value = 0.097605329781913
$ws.Cells.Item(187, 23).Value = 0.062222571450636
$ws.Cells.Item(187, 24).Value = 0.051870976887232
$ws.Cells.Item(187, 25).Value = 0.070509939810704
$ws.Cells.Item(187, 26).Value = 0.0706006328764
$ws.Cells.Item(187, 27).Value = 0.071521281839281
$ws.Cells.Item(187, 28).Value = 0.085773443327647
$ws.Cells.Item(187, 29).Value = 0.08732909313913199
$ws.Cells.Item(187, 30).Value = 0.10987856679507
$ws.Cells.Item(187, 31).Value = 0.090393774322862
$ws.Cells.Item(187, 32).Value = 0.083222062274939
$ws.Cells.Item(187, 33).Value = 0.099689077408692
$ws.Cells.Item(187, 34).Value = 0.096636298668967
$ws.Cells.Item(187, 35).Value = 0.066400267842294
$ws.Cells.Item(187, 36).Value = 0.054629096653068
$ws.Cells.Item(187, 37).Value = 0.08250195410006
$ws.Cells.Item(187, 38).Value = 0.070547159517583
$ws.Cells.Item(187, 39).Value = 0.061011652038498
$ws.Cells.Item(187, 40).Value = 0.081346891455539
$ws.Cells.Item(187, 41).Value = 0.10461213818328
$ws.Cells.Item(187, 42).Value = 0.080210118026645
$ws.Cells.Item(187, 43).Value = 0.06211595529914
$ws.Cells.Item(187, 44).Value = 0.07737248379485601
$ws.Cells.Item(187, 45).Value = 0.068718849549332
$ws.Cells.Item(187, 46).Value = 0.08967200379729399
$ws.Cells.Item(187, 47).Value = 0.096075755568854
$ws.Cells.Item(187, 48).Value = 0.09260427380627199
$ws.Cells.Item(187, 49).Value = 0.08429580991697699
$ws.Cells.Item(187, 50).Value = 0.09659842393830601
$ws.Cells.Item(187, 51).Value = 0.062137444176925
$ws.Cells.Item(187, 52).Value = 0.06511295138536601
$ws.Cells.Item(187, 53).Value = 0.065025868373043
$ws.Cells.Item(187, 54).Value = 0.072253512520964
$ws.Cells.Item(187, 55).Value = 0.074275859670593
$ws.Cells.Item(187, 56).Value = 0.073385357644732
$ws.Cells.Item(187, 57).Value = 0.09866883850569801
# Row 188
$ws.Cells.Item(188, 2).Value = 0.07704740902346401
$ws.Cells.Item(188, 3).Value = 0.10805674914637
$ws.Cells.Item(188, 4).Value = 0.098230762322482
$ws.Cells.Item(188, 6).Value = 0.062012922128249
$ws.Cells.Item(188, 7).Value = 0.074615915556594
$ws.Cells.Item(188, 8).Value = 0.088279648087456
$ws.Cells.Item(188, 9).Value = 0.07518403786165601
$ws.Cells.Item(188, 10).Value = 0.059241536603265
$ws.Cells.Item(188, 11).Value = 0.072032742926916
$ws.Cells.Item(188, 12).Value = 0.0820809614703
$ws.Cells.Item(188, 13).Value = 0.098728382149951
$ws.Cells.Item(188, 14).Value = 0.088225178690891
$ws.Cells.Item(188, 15).Value = 0.098204980034909
$ws.Cells.Item(188, 16).Value = 0.097950426693866
$ws.Cells.Item(188, 17).Value = 0.098967373664116
$ws.Cells.Item(188, 18).Value = 0.077250123797634
$ws.Cells.Item(188, 19).Value = 0.09157896474258299
$ws.Cells.Item(188, 20).Value = 0.096110911960398
$ws.Cells.Item(188, 21).Value = 0.087110011260909
$ws.Cells.Item(188, 22).Value = 0.09932845055212
$ws.Cells.Item(188, 23).Value = 0.067290874641304
$ws.Cells.Item(188, 24).Value = 0.065555760037057
$ws.Cells.Item(188, 25).Value = 0.075506681320142
$ws.Cells.Item(188, 26).Value = 0.073154971595358
$ws.Cells.Item(188, 27).Value = 0.069487002075101
$ws.Cells.Item(188, 28).Value = 0.086543673374605
$ws.Cells.Item(188, 29).Value = 0.092115108360095
$ws.Cells.Item(188, 30).Value = 0.11140639042339
$ws.Cells.Item(188, 31).Value = 0.093566344151641
$ws.Cells.Item(188, 32).Value = 0.09213831407218701
$ws.Cells.Item(188, 33).Value = 0.10282199242312
$ws.Cells.Item(188, 34).Value = 0.097280479012587
$ws.Cells.Item(188, 35).Value = 0.072834277465888
$ws.Cells.Item(188, 36).Value = 0.071066277601953
$ws.Cells.Item(188, 37).Value = 0.087425904351005
$ws.Cells.Item(188, 38).Value = 0.073295212821026
$ws.Cells.Item(188, 39).Value = 0.07078598696109301
$ws.Cells.Item(188, 40).Value = 0.08628831645228099
$ws.Cells.Item(188, 41).Value = 0.10601600872353
$ws.Cells.Item(188, 42).Value = 0.084897554563675
$ws.Cells.Item(188, 43).Value = 0.07611194395594199
$ws.Cells.Item(188, 44).Value = 0.08380122594597
$ws.Cells.Item(188, 45).Value = 0.077620935673269
$ws.Cells.Item(188, 46).Value = 0.095831333416369
$ws.Cells.Item(188, 47).Value = 0.10782443279626
$ws.Cells.Item(188, 48).Value = 0.10008232163074
$ws.Cells.Item(188, 49).Value = 0.092299860071381
$ws.Cells.Item(188, 50).Value = 0.1050809590033
$ws.Cells.Item(188, 51).Value = 0.08124445120503899
$ws.Cells.Item(188, 52).Value = 0.08363231010660099
$ws.Cells.Item(188, 53).Value = 0.07860810718231701
$ws.Cells.Item(188, 54).Value = 0.080786186483274
$ws.Cells.Item(188, 55).Value = 0.082770904934047
$ws.Cells.Item(188, 56).Value = 0.08207694148227899
$ws.Cells.Item(188, 57).Value = 0.11080436539134
